$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 (id 45): Soldier of the gates, in Purgatory
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "Soldierofthegates"
$ws.Cells.Item(46, 3).Value = "Soldier of the gates"
$ws.Cells.Item(46, 4).Value = 2
$ws.Cells.Item(46, 5).Value = "Purgatory"
$ws.Cells.Item(46, 9).Value = 1888
$ws.Cells.Item(46, 10).Value = 816

# Row 47 (id 46): Shadow, in Delusional Memories
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "Shadow"
$ws.Cells.Item(47, 3).Value = "Shadow"
$ws.Cells.Item(47, 4).Value = 2
$ws.Cells.Item(47, 5).Value = "Delusional Memories"
$ws.Cells.Item(47, 9).Value = 1824
$ws.Cells.Item(47, 10).Value = 112
